$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912": update header metadata, fix G477, append rows 478-497 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 17:51:15"
$ws1.Range("A3").Value = "Total filas: 496"
$ws1.Cells.Item(477, 7).Value = "30/12/2025"

$sheet1NewRows = @(
    ,@("", "17:51:04", "18:00", "10_OLMOS", 9, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:01", "16_SANTA ANA", 10, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:04", "17_ROMERO", 13, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:06", "23_HERNANDEZ", 15, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:11", "16_SANTA ANA", 20, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:16", "10_OLMOS", 25, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:16", "15_ABASTO", 25, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:21", "16_SANTA ANA", 30, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:21", "26_HERNANDEZ", 30, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:25", "14_ABASTO", 34, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:28", "215C_EL PATO", 37, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:32", "11X44_ETCHEVERRY", 41, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:35", "23_HERNANDEZ", 44, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:40", "15_ABASTO", 49, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:48", "14X44_ABASTO", 57, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "18:56", "10_OLMOS", 65, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "19:05", "11_ETCHEVERRY", 74, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "19:08", "23_HERNANDEZ", 77, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "19:21", "14_ABASTO", 90, "LP1912", "30/12/2025")
    ,@("", "17:51:04", "19:21", "26_HERNANDEZ", 90, "LP1912", "30/12/2025")
)

$startRow = 478
for ($i = 0; $i -lt $sheet1NewRows.Count; $i++) {
    $r = $startRow + $i
    $row = $sheet1NewRows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---- Sheet "LP1912-215": update header metadata, append row 33 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 17:51:15"
$ws2.Range("A3").Value = "Total filas: 32"

$ws2.Cells.Item(33, 1).Value = ""
$ws2.Cells.Item(33, 2).Value = "30/12/2025"
$ws2.Cells.Item(33, 3).Value = "17:51:04"
$ws2.Cells.Item(33, 4).Value = "18:28"
$ws2.Cells.Item(33, 5).Value = "215C_EL PATO"
$ws2.Cells.Item(33, 6).Value = 37
$ws2.Cells.Item(33, 7).Value = "LP1912"

# ---- Sheet "6203-6173": update header metadata, append rows 64-66 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 30/12/2025 17:51:15"
$ws3.Range("A3").Value = "Total filas: 65"

$sheet3NewRows = @(
    ,@("", "30/12/2025", "17:51:10", "18:04", "215C_LA PLATA", 13, "L6203")
    ,@("", "30/12/2025", "17:51:15", "18:52", "215A_LA PLATA", 61, "L6173")
    ,@("", "30/12/2025", "17:51:15", "19:05", "215B_LP-P MOR-1 Y 57", 74, "L6173")
)

$startRow3 = 64
for ($i = 0; $i -lt $sheet3NewRows.Count; $i++) {
    $r = $startRow3 + $i
    $row = $sheet3NewRows[$i]
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}

